$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 76886
$ws.Range("E2").Value = 4530
$ws.Range("F2").Value = 4530
$ws.Range("G2").Value = -752
$ws.Range("H2").Value = 240
$ws.Range("I2").Value = 420
$ws.Range("J2").Value = -180
$ws.Range("K2").Value = 119574
$ws.Range("L2").Value = 86691
$ws.Range("M2").Value = 32883
$ws.Range("N2").Value = 27991
$ws.Range("O2").Value = 4892
$ws.Range("P2").Value = 10373
$ws.Range("Q2").Value = 2313
$ws.Range("R2").Value = -4285
$ws.Range("S2").Value = 2814
$ws.Range("T2").Value = 2148
$ws.Range("U2").Value = 165
$ws.Range("V2").Value = 60855
$ws.Range("W2").Value = 5.89
$ws.Range("X2").Value = 0.31
$ws.Range("Y2").Value = 1.45
$ws.Range("Z2").Value = 0.2
$ws.Range("AA2").Value = 263.64
$ws.Range("AB2").Value = 148.08
$ws.Range("AC2").Value = 202
$ws.Range("AD2").Value = 48.05
$ws.Range("AE2").Value = 13493
$ws.Range("AF2").Value = 0.72
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 207455314

# Row 3
$ws.Range("D3").Value = 59649
$ws.Range("E3").Value = -951
$ws.Range("F3").Value = 274
$ws.Range("G3").Value = -7811
$ws.Range("H3").Value = -8595
$ws.Range("I3").Value = -8191
$ws.Range("J3").Value = -404
$ws.Range("K3").Value = 113832
$ws.Range("L3").Value = 82802
$ws.Range("M3").Value = 31030
$ws.Range("N3").Value = 20214
$ws.Range("O3").Value = 10816
$ws.Range("P3").Value = 10373
$ws.Range("Q3").Value = 2367
$ws.Range("R3").Value = -969
$ws.Range("S3").Value = 820
$ws.Range("T3").Value = 2148
$ws.Range("U3").Value = 219
$ws.Range("V3").Value = 60217
$ws.Range("W3").Value = -1.59
$ws.Range("X3").Value = -14.41
$ws.Range("Y3").Value = -33.98
$ws.Range("Z3").Value = -7.36
$ws.Range("AA3").Value = 266.85
$ws.Range("AB3").Value = 69.5
$ws.Range("AC3").Value = -3948
$ws.Range("AD3").Value = -1.19
$ws.Range("AE3").Value = 9744
$ws.Range("AF3").Value = 0.48
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 207455314

# Row 4
$ws.Range("D4").Value = 57296
$ws.Range("E4").Value = 4908
$ws.Range("F4").Value = 4908
$ws.Range("G4").Value = 820
$ws.Range("H4").Value = 1160
$ws.Range("I4").Value = 632
$ws.Range("J4").Value = 528
$ws.Range("K4").Value = 100268
$ws.Range("L4").Value = 65784
$ws.Range("M4").Value = 34484
$ws.Range("N4").Value = 20605
$ws.Range("O4").Value = 13879
$ws.Range("P4").Value = 10373
$ws.Range("Q4").Value = 5130
$ws.Range("R4").Value = 9096
$ws.Range("S4").Value = -14544
$ws.Range("T4").Value = 1739
$ws.Range("U4").Value = 3390
$ws.Range("V4").Value = 44232
$ws.Range("W4").Value = 8.57
$ws.Range("X4").Value = 2.02
$ws.Range("Y4").Value = 3.1
$ws.Range("Z4").Value = 1.08
$ws.Range("AA4").Value = 190.77
$ws.Range("AB4").Value = 79.8
$ws.Range("AC4").Value = 305
$ws.Range("AD4").Value = 28.96
$ws.Range("AE4").Value = 9932
$ws.Range("AF4").Value = 0.89
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 207455314

# Row 5
$ws.Range("D5").Value = 65679
$ws.Range("E5").Value = 6608
$ws.Range("F5").Value = 6608
$ws.Range("G5").Value = 4896
$ws.Range("H5").Value = 2966
$ws.Range("I5").Value = 1486
$ws.Range("J5").Value = 1480
$ws.Range("K5").Value = 102761
$ws.Range("L5").Value = 71029
$ws.Range("M5").Value = 31732
$ws.Range("N5").Value = 15937
$ws.Range("O5").Value = 15794
$ws.Range("P5").Value = 10400
$ws.Range("Q5").Value = 6657
$ws.Range("R5").Value = -2563
$ws.Range("S5").Value = 323
$ws.Range("T5").Value = 1342
$ws.Range("U5").Value = 5315
$ws.Range("V5").Value = 46490
$ws.Range("W5").Value = 10.06
$ws.Range("X5").Value = 4.52
$ws.Range("Y5").Value = 8.130000000000001
$ws.Range("Z5").Value = 2.92
$ws.Range("AA5").Value = 223.84
$ws.Range("AB5").Value = 95.79000000000001
$ws.Range("AC5").Value = 716
$ws.Range("AD5").Value = 12.14
$ws.Range("AE5").Value = 7662
$ws.Range("AF5").Value = 1.13
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 208000119

# Row 6
$ws.Range("D6").Value = 77301
$ws.Range("E6").Value = 8481
$ws.Range("F6").Value = 8481
$ws.Range("G6").Value = 6164
$ws.Range("H6").Value = 3942
$ws.Range("I6").Value = 2464
$ws.Range("K6").Value = 110292
$ws.Range("L6").Value = 72084
$ws.Range("M6").Value = 38208
$ws.Range("N6").Value = 19169
$ws.Range("P6").Value = 10408
$ws.Range("Q6").Value = 8378
$ws.Range("R6").Value = -3128
$ws.Range("S6").Value = -4226
$ws.Range("T6").Value = 1508
$ws.Range("U6").Value = 6870
$ws.Range("V6").Value = 42688
$ws.Range("W6").Value = 10.97
$ws.Range("X6").Value = 5.1
$ws.Range("Y6").Value = 14.04
$ws.Range("Z6").Value = 3.7
$ws.Range("AA6").Value = 188.66
$ws.Range("AB6").Value = 121.27
$ws.Range("AC6").Value = 1184
$ws.Range("AD6").Value = 6.41
$ws.Range("AE6").Value = 9209
$ws.Range("AF6").Value = 0.82
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 208158077
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Row 7
$ws.Range("D7").Value = 80998
$ws.Range("E7").Value = 8451
$ws.Range("G7").Value = 6229
$ws.Range("H7").Value = 4052
$ws.Range("I7").Value = 2316
$ws.Range("K7").Value = 115838
$ws.Range("L7").Value = 73025
$ws.Range("M7").Value = 42812
$ws.Range("N7").Value = 22265
$ws.Range("P7").Value = 10410
$ws.Range("Q7").Value = 4170
$ws.Range("R7").Value = -3495
$ws.Range("S7").Value = -285
$ws.Range("T7").Value = 1667
$ws.Range("U7").Value = 1387
$ws.Range("W7").Value = 10.43
$ws.Range("X7").Value = 5
$ws.Range("Y7").Value = 11.18
$ws.Range("Z7").Value = 3.58
$ws.Range("AA7").Value = 170.57
$ws.Range("AC7").Value = 1113
$ws.Range("AD7").Value = 4.42
$ws.Range("AE7").Value = 10696
$ws.Range("AF7").Value = 0.46
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 82513
$ws.Range("E8").Value = 8371
$ws.Range("G8").Value = 6546
$ws.Range("H8").Value = 4245
$ws.Range("I8").Value = 2786
$ws.Range("K8").Value = 121378
$ws.Range("L8").Value = 72695
$ws.Range("M8").Value = 48682
$ws.Range("N8").Value = 26305
$ws.Range("P8").Value = 11370
$ws.Range("Q8").Value = 6470
$ws.Range("R8").Value = -2008
$ws.Range("S8").Value = 460
$ws.Range("T8").Value = 1533
$ws.Range("U8").Value = 4570
$ws.Range("W8").Value = 10.14
$ws.Range("X8").Value = 5.15
$ws.Range("Y8").Value = 11.47
$ws.Range("Z8").Value = 3.58
$ws.Range("AA8").Value = 149.32
$ws.Range("AC8").Value = 1338
$ws.Range("AD8").Value = 3.68
$ws.Range("AE8").Value = 12637
$ws.Range("AF8").Value = 0.39
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").Value = 83976
$ws.Range("E9").Value = 8556
$ws.Range("G9").Value = 6978
$ws.Range("H9").Value = 4532
$ws.Range("I9").Value = 2866
$ws.Range("K9").Value = 123992
$ws.Range("L9").Value = 70478
$ws.Range("M9").Value = 53515
$ws.Range("N9").Value = 29178
$ws.Range("P9").Value = 11370
$ws.Range("Q9").Value = 7798
$ws.Range("R9").Value = -2035
$ws.Range("S9").Value = -2672
$ws.Range("T9").Value = 1583
$ws.Range("U9").Value = 5327
$ws.Range("W9").Value = 10.19
$ws.Range("X9").Value = 5.4
$ws.Range("Y9").Value = 10.33
$ws.Range("Z9").Value = 3.69
$ws.Range("AA9").Value = 131.7
$ws.Range("AC9").Value = 1377
$ws.Range("AD9").Value = 3.57
$ws.Range("AE9").Value = 14017
$ws.Range("AF9").Value = 0.35
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0
$ws.Range("AI9").ClearContents()
